# Updates cryptos list values (price/volume columns) per the "Updated cryptos list" commit.
# Prices in column D that read as plain decimal numbers are written with a leading
# apostrophe so Excel stores them as literal text (matching the workbook's original
# inline-string cells, e.g. "1.00" / "0.0000370") instead of silently coercing them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.748.59"
$ws.Range("E2").Value = "  +0.71%  "

$ws.Range("D3").Value = "3.828.47"
$ws.Range("E3").Value = "  +3.98%  "

$ws.Range("E4").Value = "  -0.12%  "

$ws.Range("D5").Value = "'411.57"
$ws.Range("E5").Value = "  -1.64%  "

$ws.Range("D6").Value = "'131.06"
$ws.Range("E6").Value = "  +1.02%  "

$ws.Range("D7").Value = "3.819.85"
$ws.Range("E7").Value = "  +3.89%  "

$ws.Range("D8").Value = "'0.614"
$ws.Range("E8").Value = "  -4.02%  "

$ws.Range("D9").Value = "'1.00"
$ws.Range("E9").Value = "  +0.01%  "

$ws.Range("D10").Value = "'0.733"
$ws.Range("E10").Value = "  -4.21%  "

$ws.Range("E11").Value = "  -4.86%  "

$ws.Range("D12").Value = "'0.0000370"
$ws.Range("E12").Value = "  -4.68%  "

$ws.Range("D13").Value = "'41.09"
$ws.Range("E13").Value = "  -4.32%  "

$ws.Range("D14").Value = "4.440.59"
$ws.Range("E14").Value = "  +3.81%  "

$ws.Range("D15").Value = "'10.01"
$ws.Range("E15").Value = "  -5.42%  "

$ws.Range("D16").Value = "'15.45"
$ws.Range("E16").Value = "  +16.82%  "

$ws.Range("E17").Value = "  -0.99%  "

$ws.Range("D18").Value = "3.826.54"
$ws.Range("E18").Value = "  +4.03%  "

$ws.Range("D19").Value = "'19.69"
$ws.Range("E19").Value = "  -3.93%  "

$ws.Range("D20").Value = "67.226.40"
$ws.Range("E20").Value = "  +1.17%  "

$ws.Range("E21").Value = "  -4.58%  "

$ws.Range("D22").Value = "'412.42"
$ws.Range("E22").Value = "  -6.80%  "

$ws.Range("D23").Value = "'14.45"
$ws.Range("E23").Value = "  -11.78%  "

$ws.Range("D24").Value = "'85.31"
$ws.Range("E24").Value = "  -4.99%  "

$ws.Range("E25").Value = "  -2.54%  "

$ws.Range("D26").Value = "'5.77"
$ws.Range("E26").Value = "  +13.80%  "

$ws.Range("D27").Value = "'36.51"
$ws.Range("E27").Value = "  -1.74%  "

$ws.Range("D28").Value = "'3.13"
$ws.Range("E28").Value = "  -5.93%  "

$ws.Range("D29").Value = "'9.42"
$ws.Range("E29").Value = "  -7.49%  "

$ws.Range("D30").Value = "'679.93"
$ws.Range("E30").Value = "  +4.50%  "

$ws.Range("E31").Value = "  -2.56%  "

$ws.Range("D32").Value = "'12.45"
$ws.Range("E32").Value = "  -1.79%  "

$ws.Range("E33").Value = "  -1.79%  "

$ws.Range("D34").Value = "'7.18"
$ws.Range("E34").Value = "  -1.47%  "

$ws.Range("E35").Value = "  -8.35%  "

$ws.Range("D36").Value = "'38.66"
$ws.Range("E36").Value = "  -6.99%  "

$ws.Range("E37").Value = "  -0.01%  "

$ws.Range("E38").Value = "  +9.54%  "

$ws.Range("D39").Value = "'54.99"
$ws.Range("E39").Value = "  -3.96%  "

$ws.Range("D40").Value = "'3.09"
$ws.Range("E40").Value = "  -0.64%  "

$ws.Range("E41").Value = "  -6.67%  "

$ws.Range("D42").Value = "'0.997"
$ws.Range("E42").Value = "  -0.12%  "

$ws.Range("E43").Value = "  -8.69%  "

$ws.Range("D44").Value = "'149.76"
$ws.Range("E44").Value = "  +0.73%  "

$ws.Range("D45").Value = "'4.49"
$ws.Range("E45").Value = "  +2.44%  "

$ws.Range("E46").Value = "  -2.83%  "

$ws.Range("D47").Value = "'3.14"
$ws.Range("E47").Value = "  +13.85%  "

$ws.Range("B48").Value = "ARBITRUM"
$ws.Range("C48").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D48").Value = "'2.08"
$ws.Range("E48").Value = "  -0.84%  "

$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "'26.49"
$ws.Range("E49").Value = "  -9.68%  "

$ws.Range("D50").Value = "'2.79"

$ws.Range("D51").Value = "'2.56"
$ws.Range("E51").Value = "  -3.87%  "
